# StructureDefinition-industry-classification.xlsx update
# - Bump Version 5.0.0 -> 6.0.0
# - Bump Date to new publish timestamp
# - Replace empty "Contact" rows with real Publisher ("Alvearie Team") and a new
#   "Jurisdiction" ("United States of America") row
# - Update the root Extension's Short/Definition on the Elements sheet to the
#   real title/description instead of the generic placeholders

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Row 11 was a duplicate "Contact / No display for ContactDetail" placeholder
# row; remove it so everything below shifts up by one (A1:B21 -> A1:B20).
$meta.Rows.Item(11).Delete()

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher (row 9 used to be "Publisher" / blank)
$meta.Range("B9").Value = "Alvearie Team"

# Jurisdiction (row 10 used to be the first "Contact" placeholder row)
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition get the real title & description
# instead of the generic "Extension" / "An Extension" placeholders.
$elements.Range("K2").Value = "Industry Classification"
$elements.Range("L2").Value = "Industry classification code which can be based on the North American Industry Classification System (NAICS)"
